$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 154
$ws.Cells.Item(154, 2).Value2 = 64350
$ws.Cells.Item(154, 3).Value2 = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Cells.Item(154, 4).Value2 = 66.44
$ws.Cells.Item(154, 5).Value2 = 70.63
$ws.Cells.Item(154, 6).Value2 = 101
$ws.Cells.Item(154, 7).Value2 = 6710.44

# Row 155
$ws.Cells.Item(155, 2).Value2 = 57756
$ws.Cells.Item(155, 3).Value2 = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Cells.Item(155, 4).Value2 = 66.44
$ws.Cells.Item(155, 5).Value2 = 79.37
$ws.Cells.Item(155, 6).Value2 = -100
$ws.Cells.Item(155, 7).Value2 = -6644

# Row 156
$ws.Cells.Item(156, 2).Value2 = 53925
$ws.Cells.Item(156, 3).Value2 = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Cells.Item(156, 4).Value2 = 66.44
$ws.Cells.Item(156, 5).Value2 = 79.37
$ws.Cells.Item(156, 6).Value2 = 1
$ws.Cells.Item(156, 7).Value2 = 66.44

# Row 176
$ws.Cells.Item(176, 2).Value2 = 57552
$ws.Cells.Item(176, 3).Value2 = "DAB-Real Activ Coconut Water Tetra 1000ml"
$ws.Cells.Item(176, 4).Value2 = 120.69
$ws.Cells.Item(176, 5).Value2 = 136.86
$ws.Cells.Item(176, 6).Value2 = -5
$ws.Cells.Item(176, 7).Value2 = -603.45

# Row 177
$ws.Cells.Item(177, 2).Value2 = 64329
$ws.Cells.Item(177, 3).Value2 = "DAB-Real Activ Coconut Water Tetra 1000ml"
$ws.Cells.Item(177, 4).Value2 = 120.69
$ws.Cells.Item(177, 5).Value2 = 128.32
$ws.Cells.Item(177, 6).Value2 = 6
$ws.Cells.Item(177, 7).Value2 = 724.14

# Row 256
$ws.Cells.Item(256, 2).Value2 = 64979
$ws.Cells.Item(256, 3).Value2 = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Cells.Item(256, 4).Value2 = 295.75
$ws.Cells.Item(256, 5).Value2 = 314.41
$ws.Cells.Item(256, 6).Value2 = 82
$ws.Cells.Item(256, 7).Value2 = 24251.5

# Row 257
$ws.Cells.Item(257, 2).Value2 = 48719
$ws.Cells.Item(257, 3).Value2 = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Cells.Item(257, 4).Value2 = 295.75
$ws.Cells.Item(257, 5).Value2 = 353.35
$ws.Cells.Item(257, 6).Value2 = -81
$ws.Cells.Item(257, 7).Value2 = -23955.75

# Row 305
$ws.Cells.Item(305, 2).Value2 = 57854
$ws.Cells.Item(305, 3).Value2 = "HUL-3Roses Dust [C] 500G Relaunch"
$ws.Cells.Item(305, 4).Value2 = 305.84
$ws.Cells.Item(305, 5).Value2 = 325.16
$ws.Cells.Item(305, 6).Value2 = 2
$ws.Cells.Item(305, 7).Value2 = 611.6799999999999

# Row 306
$ws.Cells.Item(306, 2).Value2 = 62997
$ws.Cells.Item(306, 3).Value2 = "HUL-3Roses Dust [C] 500G Relaunch"
$ws.Cells.Item(306, 4).Value2 = 305.84
$ws.Cells.Item(306, 5).Value2 = 325.16
$ws.Cells.Item(306, 6).Value2 = 72
$ws.Cells.Item(306, 7).Value2 = 22020.48

# Row 309
$ws.Cells.Item(309, 2).Value2 = 61610
$ws.Cells.Item(309, 3).Value2 = "HUL-Bru Inst Poly 50g"
$ws.Cells.Item(309, 4).Value2 = 102.71
$ws.Cells.Item(309, 5).Value2 = 122.71
$ws.Cells.Item(309, 6).Value2 = -58
$ws.Cells.Item(309, 7).Value2 = -5957.18

# Row 310
$ws.Cells.Item(310, 2).Value2 = 63565
$ws.Cells.Item(310, 3).Value2 = "HUL-Bru Inst Poly 50g"
$ws.Cells.Item(310, 4).Value2 = 102.71
$ws.Cells.Item(310, 5).Value2 = 109.19
$ws.Cells.Item(310, 6).Value2 = 60
$ws.Cells.Item(310, 7).Value2 = 6162.6

# Row 342
$ws.Cells.Item(342, 2).Value2 = 63531
$ws.Cells.Item(342, 3).Value2 = "HUL-Kissan Pineapple Jam 500G"
$ws.Cells.Item(342, 4).Value2 = 143.48
$ws.Cells.Item(342, 5).Value2 = 152.53
$ws.Cells.Item(342, 6).Value2 = 80
$ws.Cells.Item(342, 7).Value2 = 11478.4

# Row 343
$ws.Cells.Item(343, 2).Value2 = 57802
$ws.Cells.Item(343, 3).Value2 = "HUL-Kissan Pineapple Jam 500G"
$ws.Cells.Item(343, 4).Value2 = 143.48
$ws.Cells.Item(343, 5).Value2 = 162.71
$ws.Cells.Item(343, 6).Value2 = -79
$ws.Cells.Item(343, 7).Value2 = -11334.92

# Row 344
$ws.Cells.Item(344, 2).Value2 = 63571
$ws.Cells.Item(344, 3).Value2 = "HUL-Kissan Pineapple Jam 500G"
$ws.Cells.Item(344, 4).Value2 = 143.48
$ws.Cells.Item(344, 5).Value2 = 152.53
$ws.Cells.Item(344, 6).Value2 = 29
$ws.Cells.Item(344, 7).Value2 = 4160.92

# Row 371
$ws.Cells.Item(371, 2).Value2 = 61608
$ws.Cells.Item(371, 3).Value2 = "HUL-Lux Radiant Glow 3*150g"
$ws.Cells.Item(371, 4).Value2 = 129.01
$ws.Cells.Item(371, 5).Value2 = 154.12
$ws.Cells.Item(371, 6).Value2 = -56
$ws.Cells.Item(371, 7).Value2 = -7224.56

# Row 372
$ws.Cells.Item(372, 2).Value2 = 63564
$ws.Cells.Item(372, 3).Value2 = "HUL-Lux Radiant Glow 3*150g"
$ws.Cells.Item(372, 4).Value2 = 129.01
$ws.Cells.Item(372, 5).Value2 = 137.16
$ws.Cells.Item(372, 6).Value2 = 57
$ws.Cells.Item(372, 7).Value2 = 7353.57

# Row 381
$ws.Cells.Item(381, 2).Value2 = 62865
$ws.Cells.Item(381, 3).Value2 = "HUL-Rap Refresh Bolt 1Kg"
$ws.Cells.Item(381, 4).Value2 = 79.81
$ws.Cells.Item(381, 5).Value2 = 95.34999999999999
$ws.Cells.Item(381, 6).Value2 = 151
$ws.Cells.Item(381, 7).Value2 = 12051.31

# Row 382
$ws.Cells.Item(382, 2).Value2 = 57817
$ws.Cells.Item(382, 3).Value2 = "HUL-Rap Refresh Bolt 1Kg"
$ws.Cells.Item(382, 4).Value2 = 79.81
$ws.Cells.Item(382, 5).Value2 = 95.34999999999999
$ws.Cells.Item(382, 6).Value2 = 3
$ws.Cells.Item(382, 7).Value2 = 239.43

# Row 411
$ws.Cells.Item(411, 2).Value2 = 63007
$ws.Cells.Item(411, 3).Value2 = "HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp"
$ws.Cells.Item(411, 4).Value2 = 171.33
$ws.Cells.Item(411, 5).Value2 = 204.69
$ws.Cells.Item(411, 6).Value2 = 984
$ws.Cells.Item(411, 7).Value2 = 168588.72

# Row 412
$ws.Cells.Item(412, 2).Value2 = 57856
$ws.Cells.Item(412, 3).Value2 = "HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp"
$ws.Cells.Item(412, 4).Value2 = 171.33
$ws.Cells.Item(412, 5).Value2 = 204.69
$ws.Cells.Item(412, 6).Value2 = 2
$ws.Cells.Item(412, 7).Value2 = 342.66

# Row 423
$ws.Cells.Item(423, 2).Value2 = 63102
$ws.Cells.Item(423, 3).Value2 = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Cells.Item(423, 4).Value2 = 59.47
$ws.Cells.Item(423, 5).Value2 = 71.05
$ws.Cells.Item(423, 6).Value2 = 36
$ws.Cells.Item(423, 7).Value2 = 2140.92

# Row 424
$ws.Cells.Item(424, 2).Value2 = 53082
$ws.Cells.Item(424, 3).Value2 = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Cells.Item(424, 4).Value2 = 59.47
$ws.Cells.Item(424, 5).Value2 = 71.05
$ws.Cells.Item(424, 6).Value2 = 1
$ws.Cells.Item(424, 7).Value2 = 59.47

# Row 528
$ws.Cells.Item(528, 2).Value2 = 58047
$ws.Cells.Item(528, 3).Value2 = "KUS-Floor Wiper"
$ws.Cells.Item(528, 4).Value2 = 105.54
$ws.Cells.Item(528, 5).Value2 = 126.1
$ws.Cells.Item(528, 6).Value2 = 54
$ws.Cells.Item(528, 7).Value2 = 5699.16

# Row 529
$ws.Cells.Item(529, 2).Value2 = 47097
$ws.Cells.Item(529, 3).Value2 = "KUS-Floor Wiper"
$ws.Cells.Item(529, 4).Value2 = 112.28
$ws.Cells.Item(529, 5).Value2 = 134.16
$ws.Cells.Item(529, 6).Value2 = 15
$ws.Cells.Item(529, 7).Value2 = 1684.2

# Row 578
$ws.Cells.Item(578, 2).Value2 = 64915
$ws.Cells.Item(578, 3).Value2 = "CRE-Cremica Chocolate Cream 150Gm"
$ws.Cells.Item(578, 4).Value2 = 19.73
$ws.Cells.Item(578, 5).Value2 = 20.98
$ws.Cells.Item(578, 6).Value2 = 40
$ws.Cells.Item(578, 7).Value2 = 789.2

# Row 579
$ws.Cells.Item(579, 2).Value2 = 45695
$ws.Cells.Item(579, 3).Value2 = "CRE-Cremica Chocolate Cream 150Gm"
$ws.Cells.Item(579, 4).Value2 = 19.73
$ws.Cells.Item(579, 5).Value2 = 23.58
$ws.Cells.Item(579, 6).Value2 = -36
$ws.Cells.Item(579, 7).Value2 = -710.28

# Row 585
$ws.Cells.Item(585, 2).Value2 = 45718
$ws.Cells.Item(585, 3).Value2 = "CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm"
$ws.Cells.Item(585, 4).Value2 = 16.22
$ws.Cells.Item(585, 5).Value2 = 19.38
$ws.Cells.Item(585, 6).Value2 = -294
$ws.Cells.Item(585, 7).Value2 = -4768.68

# Row 586
$ws.Cells.Item(586, 2).Value2 = 64927
$ws.Cells.Item(586, 3).Value2 = "CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm"
$ws.Cells.Item(586, 4).Value2 = 16.22
$ws.Cells.Item(586, 5).Value2 = 17.26
$ws.Cells.Item(586, 6).Value2 = 295
$ws.Cells.Item(586, 7).Value2 = 4784.9

# Row 593
$ws.Cells.Item(593, 2).Value2 = 64919
$ws.Cells.Item(593, 3).Value2 = "CRE-Cremica Pista Almond Cookies (75 +25Gm)"
$ws.Cells.Item(593, 4).Value2 = 26.3
$ws.Cells.Item(593, 5).Value2 = 27.97
$ws.Cells.Item(593, 6).Value2 = 224
$ws.Cells.Item(593, 7).Value2 = 5891.2

# Row 594
$ws.Cells.Item(594, 2).Value2 = 45702
$ws.Cells.Item(594, 3).Value2 = "CRE-Cremica Pista Almond Cookies (75 +25Gm)"
$ws.Cells.Item(594, 4).Value2 = 26.3
$ws.Cells.Item(594, 5).Value2 = 31.43
$ws.Cells.Item(594, 6).Value2 = -215
$ws.Cells.Item(594, 7).Value2 = -5654.5

# Row 679
$ws.Cells.Item(679, 2).Value2 = 64810
$ws.Cells.Item(679, 3).Value2 = "PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)"
$ws.Cells.Item(679, 4).Value2 = 273.92
$ws.Cells.Item(679, 5).Value2 = 291.22
$ws.Cells.Item(679, 6).Value2 = 7
$ws.Cells.Item(679, 7).Value2 = 1917.44

# Row 680
$ws.Cells.Item(680, 2).Value2 = 53319
$ws.Cells.Item(680, 3).Value2 = "PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)"
$ws.Cells.Item(680, 4).Value2 = 273.92
$ws.Cells.Item(680, 5).Value2 = 310.64
$ws.Cells.Item(680, 6).Value2 = -6
$ws.Cells.Item(680, 7).Value2 = -1643.52

# Row 712
$ws.Cells.Item(712, 2).Value2 = 64830
$ws.Cells.Item(712, 3).Value2 = "Rasna Nagpur Orange (32 Glass)"
$ws.Cells.Item(712, 4).Value2 = 32.83
$ws.Cells.Item(712, 5).Value2 = 34.9
$ws.Cells.Item(712, 6).Value2 = 117
$ws.Cells.Item(712, 7).Value2 = 3841.11

# Row 713
$ws.Cells.Item(713, 2).Value2 = 60022
$ws.Cells.Item(713, 3).Value2 = "Rasna Nagpur Orange (32 Glass)"
$ws.Cells.Item(713, 4).Value2 = 32.83
$ws.Cells.Item(713, 5).Value2 = 37.22
$ws.Cells.Item(713, 6).Value2 = -113
$ws.Cells.Item(713, 7).Value2 = -3709.79

# Row 864
$ws.Cells.Item(864, 2).Value2 = 65079
$ws.Cells.Item(864, 3).Value2 = "Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm"
$ws.Cells.Item(864, 4).Value2 = 40.87
$ws.Cells.Item(864, 5).Value2 = 43.44
$ws.Cells.Item(864, 6).Value2 = 21
$ws.Cells.Item(864, 7).Value2 = 858.27

# Row 865
$ws.Cells.Item(865, 2).Value2 = 54751
$ws.Cells.Item(865, 3).Value2 = "Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm"
$ws.Cells.Item(865, 4).Value2 = 40.87
$ws.Cells.Item(865, 5).Value2 = 46.34
$ws.Cells.Item(865, 6).Value2 = -19
$ws.Cells.Item(865, 7).Value2 = -776.53
